$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.902.15"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.40"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7415"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.39"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3144"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07205"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08321"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7494"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.886.63"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.378"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.17"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.115"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.913.80"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "247.45"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.54"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007835"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.127.78"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.998"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1552"
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.290"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.00"
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.64"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.023"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.483"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.577"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.533"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.227"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05326"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7495"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01964"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.753"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4547"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.131.30"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.126"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.35"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8589"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.45"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.617"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.506"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.028.17"
$ws.Range("E51").Value = "  +0.48%  "
